$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 8029835
$ws.Range("I33").Value = 2813825.8
$ws.Range("K33").Value = 2813825.8
$ws.Range("M33").Value = -2813596.8
$ws.Range("H62").Value = 7746.25
$ws.Range("I62").Value = 6995
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 6995
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -6371
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 7746.25
$ws.Range("I65").Value = 6995
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 34975
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -31855
$ws.Range("N65").Value = -56240
$ws.Range("H98").Value = 2160.475
$ws.Range("I98").Value = 1880.697
$ws.Range("J98").Value = 3479.4285
$ws.Range("K98").Value = 1880.697
$ws.Range("L98").Value = 3479.4285
$ws.Range("M98").Value = -382.6969999999999
$ws.Range("N98").Value = -6475.4285
$ws.Range("H103").Value = 427.2857
$ws.Range("I103").Value = 425
$ws.Range("J103").Value = 433
$ws.Range("K103").Value = 1275
$ws.Range("L103").Value = 1299
$ws.Range("M103").Value = -689
$ws.Range("N103").Value = -2471
$ws.Range("H122").Value = 2160.475
$ws.Range("I122").Value = 1880.697
$ws.Range("J122").Value = 3479.4285
$ws.Range("K122").Value = 5642.090999999999
$ws.Range("L122").Value = 10438.2855
$ws.Range("M122").Value = -3192.090999999999
$ws.Range("N122").Value = -15338.2855
$ws.Range("H135").Value = 1173.6
$ws.Range("I135").Value = 1173.6
$ws.Range("K135").Value = 10562.4
$ws.Range("M135").Value = -8027.4
$ws.Range("H137").Value = 1305006.1
$ws.Range("I137").Value = 6179.7812
$ws.Range("J137").Value = 2228616
$ws.Range("K137").Value = 18539.3436
$ws.Range("L137").Value = 6685848
$ws.Range("M137").Value = -15989.3436
$ws.Range("N137").Value = -6690948
$ws.Range("H138").Value = 12443.385
$ws.Range("I138").Value = 25089.8
$ws.Range("J138").Value = 4539.375
$ws.Range("K138").Value = 75269.39999999999
$ws.Range("L138").Value = 13618.125
$ws.Range("M138").Value = -70129.39999999999
$ws.Range("N138").Value = -23898.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2277.4302
$ws.Range("I32").Value = 1205.5343
$ws.Range("J32").Value = 8296.538
$ws.Range("K32").Value = 1205.5343
$ws.Range("L32").Value = 8296.538
$ws.Range("M32").Value = -918.5343
$ws.Range("N32").Value = -8870.538
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = $null
$ws.Range("H97").Value = 5295.5
$ws.Range("I97").Value = 6059.421
$ws.Range("K97").Value = 6059.421
$ws.Range("M97").Value = -5563.421
$ws.Range("H132").Value = 2668.9714
$ws.Range("I132").Value = 1610.9667
$ws.Range("J132").Value = 9017
$ws.Range("K132").Value = 4832.9001
$ws.Range("L132").Value = 27051
$ws.Range("M132").Value = -2302.9001
$ws.Range("N132").Value = -32111

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6761.2964
$ws.Range("I99").Value = 8097.9
$ws.Range("J99").Value = 2942.4285
$ws.Range("K99").Value = 8097.9
$ws.Range("L99").Value = 2942.4285
$ws.Range("M99").Value = -6599.9
$ws.Range("N99").Value = -5938.4285
$ws.Range("H122").Value = 51965.918
$ws.Range("J122").Value = 51965.918
$ws.Range("L122").Value = 51965.918
$ws.Range("N122").Value = -61765.918

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4589.3403
$ws.Range("I31").Value = 2922.45
$ws.Range("K31").Value = 2922.45
$ws.Range("M31").Value = -2627.45
$ws.Range("H34").Value = 4589.3403
$ws.Range("I34").Value = 2922.45
$ws.Range("K34").Value = 2922.45
$ws.Range("M34").Value = -2720.45
$ws.Range("H94").Value = 1322.2307
$ws.Range("I94").Value = 1178.8
$ws.Range("J94").Value = 1411.875
$ws.Range("K94").Value = 1178.8
$ws.Range("L94").Value = 1411.875
$ws.Range("M94").Value = -727.8
$ws.Range("N94").Value = -2313.875
$ws.Range("H105").Value = 2994.8667
$ws.Range("I105").Value = 2401.182
$ws.Range("K105").Value = 2401.182
$ws.Range("M105").Value = -654.1819999999998
$ws.Range("H132").Value = 3540.7
$ws.Range("I132").Value = 4057.2307
$ws.Range("J132").Value = 2581.4285
$ws.Range("K132").Value = 12171.6921
$ws.Range("L132").Value = 7744.2855
$ws.Range("M132").Value = -9641.6921
$ws.Range("N132").Value = -12804.2855
$ws.Range("H134").Value = 1750
$ws.Range("I134").Value = 1833.3334
$ws.Range("K134").Value = 5500.0002
$ws.Range("M134").Value = -2965.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2809.3125
$ws.Range("J68").Value = 3989.8
$ws.Range("L68").Value = 11969.4
$ws.Range("N68").Value = -13591.4
$ws.Range("H71").Value = 2809.3125
$ws.Range("J71").Value = 3989.8
$ws.Range("L71").Value = 35908.2
$ws.Range("N71").Value = -44020.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 440.47827
$ws.Range("I97").Value = 465.6
$ws.Range("J97").Value = 393.375
$ws.Range("K97").Value = 465.6
$ws.Range("L97").Value = 393.375
$ws.Range("M97").Value = 30.39999999999998
$ws.Range("N97").Value = -1385.375
$ws.Range("H141").Value = 79949.2
$ws.Range("I141").Value = 89000
$ws.Range("J141").Value = 77686.5
$ws.Range("K141").Value = 89000
$ws.Range("L141").Value = 77686.5
$ws.Range("M141").Value = -83820
$ws.Range("N141").Value = -88046.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5381.357
$ws.Range("I68").Value = 5659.778
$ws.Range("J68").Value = 4880.2
$ws.Range("K68").Value = 5659.778
$ws.Range("L68").Value = 4880.2
$ws.Range("M68").Value = -4910.778
$ws.Range("N68").Value = -6378.2
$ws.Range("H71").Value = 5381.357
$ws.Range("I71").Value = 5659.778
$ws.Range("J71").Value = 4880.2
$ws.Range("K71").Value = 28298.89
$ws.Range("L71").Value = 24401
$ws.Range("M71").Value = -24554.89
$ws.Range("N71").Value = -31889

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 92798.75
$ws.Range("I46").Value = 22222
$ws.Range("J46").Value = 99214.82000000001
$ws.Range("K46").Value = 22222
$ws.Range("L46").Value = 99214.82000000001
$ws.Range("M46").Value = -21991
$ws.Range("N46").Value = -99676.82000000001
$ws.Range("H61").Value = 2200
$ws.Range("I61").Value = 2200
$ws.Range("K61").Value = 2200
$ws.Range("M61").Value = -1908
$ws.Range("H132").Value = 1460.3334
$ws.Range("I132").Value = 1297.1482
$ws.Range("J132").Value = 2929
$ws.Range("K132").Value = 3891.4446
$ws.Range("L132").Value = 8787
$ws.Range("M132").Value = -1361.4446
$ws.Range("N132").Value = -13847
$ws.Range("H134").Value = 92798.75
$ws.Range("I134").Value = 22222
$ws.Range("J134").Value = 99214.82000000001
$ws.Range("K134").Value = 66666
$ws.Range("L134").Value = 297644.46
$ws.Range("M134").Value = -64131
$ws.Range("N134").Value = -302714.46
$ws.Range("H139").Value = 162942.42
$ws.Range("I139").Value = 100650
$ws.Range("K139").Value = 100650
$ws.Range("M139").Value = -95510
